$d = $word.ActiveDocument

# =====================================================================
# Change 1: title block  "з дисципліни “" -> "н" + bookmark(_GoBack) +
#           "а тему" + " “"   (i.e. overall text becomes "на тему “")
# =====================================================================
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("з дисципліни", $true, $false, $false, $false, $false, $true, 1, $false, "на тему", 2)
Write-Output "Change1 replace: $found1"

$locate1 = $d.Content
$found1b = $locate1.Find.Execute("на тему", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Change1 locate: $found1b"

# Split "н" | "а тему "" by dropping a _GoBack bookmark right after the
# first letter (this naturally breaks the run there).
$bmPos = $locate1.Start + 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Output "Change1 bookmark added"

# Split "а тему" from " "" by toggling Bold off/on across the trailing
# two characters (space + opening quote char) -- this forces a distinct
# run without leaving any formatting residue behind.
$tailStart = $locate1.Start + 7
$tailEnd = $locate1.Start + 9
$tail = $d.Range($tailStart, $tailEnd)
$tail.Font.Bold = 0
$tail2 = $d.Range($tailStart, $tailEnd)
$tail2.Font.Bold = 1
Write-Output "Change1 tail split done"

# =====================================================================
# Change 2: second occurrence of the lab-title quote -- merge the lone
#           leading-space run into the following "Запити..." run,
#           without touching the quote-mark runs on either side.
# =====================================================================
$anchor2 = $d.Content
$foundAnchor2 = $anchor2.Find.Execute("на тему:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Change2 anchor: $foundAnchor2"

$search2 = $d.Range($anchor2.End, $d.Content.End)
$found2 = $search2.Find.Execute("Запити на вибір даних з таблиць бази даних ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Change2 locate: $found2"

# Neutralise the closing-quote run's formatting momentarily so the merge
# triggered below does not pull it in too.
$closeQuote = $d.Range($search2.End, $search2.End + 1)
$closeQuote.Font.Bold = 0
$mergeRng = $d.Range($search2.Start, $search2.End)
$foundMerge = $mergeRng.Find.Execute("Запити", $true, $false, $false, $false, $false, $true, 1, $false, "Запити", 2)
Write-Output "Change2 merge: $foundMerge"
$closeQuote2 = $d.Range($search2.End, $search2.End + 1)
$closeQuote2.Font.Bold = 1
Write-Output "Change2 restore done"

# =====================================================================
# Change 3: "Оберем" + bookmark(_GoBack) + "о" -> single run "Оберемо"
# =====================================================================
$rng3 = $d.Content
$found3a = $rng3.Find.Execute("Оберем", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Change3 locate: $found3a"

$mergeRng3 = $d.Range($rng3.Start, $rng3.End + 1)
$found3 = $mergeRng3.Find.Execute("Оберемо", $true, $false, $false, $false, $false, $true, 1, $false, "Оберемо", 2)
Write-Output "Change3 merge: $found3"

Write-Output "ALL DONE"
